# Auto commit at 2025-09-01  9:47:23.73
# Adds the 2025-08-31 (serial 45900) daily totals for the two charging
# stations as rows 62 and 63, and moves the worksheet view/selection
# down to keep the newly-added rows visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 62: 四方坪站 (Sifangping station) -------------------------------
$ws.Cells.Item(62, 1).Value = 45900
$ws.Cells.Item(62, 2).Value = "四方坪站"
$ws.Cells.Item(62, 3).Value = 10299.58
$ws.Cells.Item(62, 4).Value = 8861.17
$ws.Cells.Item(62, 5).Value = 3546.75
$ws.Cells.Item(62, 6).Value = 438

# ---- Row 63: 高岭站 (Gaoling station) -------------------------------------
$ws.Cells.Item(63, 1).Value = 45900
$ws.Cells.Item(63, 2).Value = "高岭站"
$ws.Cells.Item(63, 3).Value = 5351.64
$ws.Cells.Item(63, 4).Value = 4463.61
$ws.Cells.Item(63, 5).Value = 1393.72
$ws.Cells.Item(63, 6).Value = 171

# ---- Update the saved view/selection state -------------------------------
# Scroll so row 52 becomes the top visible row, and move the active
# selection to H58 (matching the author's on-screen state at save time).
$win = $excel.ActiveWindow
$ws.Range("H58").Select() | Out-Null
$win.ScrollRow = 52
$win.ScrollColumn = 1
